$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert two new columns before the current column D ("Waypart %
#    Util / Ways" sub-table). Everything from the old D column
#    onward (D..J) shifts right to F..L.
# ------------------------------------------------------------------
$ws.Range("D1:E1").EntireColumn.Insert()

# The insert copies formatting from the left neighbour (bold style of
# C2), but the new header cells D2/E2 ("% Util." / "Ways") are meant
# to be unstyled, so strip the inherited formatting back off.
$ws.Range("D2:E2").ClearFormats()

# ------------------------------------------------------------------
# 2. Fill in the new "% Util." / "Ways" helper columns.
# ------------------------------------------------------------------
$ws.Range("M2").Value = "EquiBlissPart"
$ws.Range("D2").Value = "% Util."
$ws.Range("E2").Value = "Ways"
$ws.Range("N2").Value = "UnevenWays"

$ws.Range("D3").Formula = '=C3/SUM($C$3:$C$6)'
$ws.Range("D4").Formula = '=C4/SUM($C$3:$C$6)'
$ws.Range("D5").Formula = '=C5/SUM($C$3:$C$6)'
$ws.Range("D6").Formula = '=C6/SUM($C$3:$C$6)'

$ws.Range("E3").Formula = '=D3*8'
$ws.Range("E4").Formula = '=D4*8'
$ws.Range("E5").Formula = '=D5*8'
$ws.Range("E6").Formula = '=D6*8'

# ------------------------------------------------------------------
# 3. Add the new "EquiBlissPart" (M) and "UnevenWays" (N) scheduler
#    columns, mirroring the layout of the other scheduler columns.
# ------------------------------------------------------------------
$ws.Range("M2:N2").Font.Bold = $true

$ws.Range("M3").Value = 1768461
$ws.Range("M4").Value = 33845279
$ws.Range("M5").Value = 31370079
$ws.Range("M6").Value = 18987613

$ws.Range("N3").Value = 1768461
$ws.Range("N4").Value = 33845279
$ws.Range("N5").Value = 31370079
$ws.Range("N6").Value = 18987613

$ws.Range("M8").Formula = '=MAX(($C$3/M3), ($C$4/M4), ($C$5/M5), ($C$6/M6))'
$ws.Range("N8").Formula = '=MAX(($C$3/N3), ($C$4/N4), ($C$5/N5), ($C$6/N6))'

$ws.Range("M9").Formula = '=(M3/$C$3)+(M4/$C$4)+(M5/$C$5)+(M6/$C$6)'
$ws.Range("N9").Formula = '=(N3/$C$3)+(N4/$C$4)+(N5/$C$5)+(N6/$C$6)'

$ws.Range("M10").Formula = '=AVERAGE(M3:M6)'
$ws.Range("N10").Formula = '=AVERAGE(N3:N6)'

$ws.Range("M12").Formula = '=M8/$H$8'
$ws.Range("N12").Formula = '=N8/$H$8'

$ws.Range("M13").Formula = '=M9/$H$9'
$ws.Range("N13").Formula = '=N9/$H$9'

$ws.Range("M14").Formula = '=$H$10/M10'
$ws.Range("N14").Formula = '=$H$10/N10'

# ------------------------------------------------------------------
# 4. Misc cosmetic bits that the diff also records: the selected
#    cell in the sheet view.
# ------------------------------------------------------------------
[void]$ws.Range("Q11").Select()
